# Generate Report for Handback
# Update the "Latest HO / Correspond Handoff / Correspond Handback" datetime
# stamps for the 2dd4ca1a... file row after a new handback round completed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 2dd4ca1a...md file; column G = "Latest HO Xliff Generate Date"
$overview.Range("G3").Value = "2016-09-04 08:52:32"

# zh-cn sheet: row 3 is the 2dd4ca1a...md file
# H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$zhcn.Range("H3").Value = "2016-09-04 08:52:28"
$zhcn.Range("K3").Value = "2016-09-04 08:52:52"

# de-de sheet: row 3 is the 2dd4ca1a...md file
# H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$dede.Range("H3").Value = "2016-09-04 08:52:32"
$dede.Range("K3").Value = "2016-09-04 08:52:59"
